# Insert a new weekly price record as row 572 in the "Acelga" price
# history sheet, pushing the existing rows 572:598 down to 573:599.
#
# The new record (Femacal de La Calera, Coquimbo, Acelga) is dated
# serial 45147 (2023-08-09) with unit price range 3000-3000-3000 and P=500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 572..598 down to 573..599, leaving a blank row 572.
$ws.Rows.Item(572).Insert()

# Populate the new row 572 with the inserted record's values.
$ws.Cells.Item(572, 1).Value = 3
$ws.Cells.Item(572, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(572, 3).Value = "Coquimbo"
$ws.Cells.Item(572, 4).Value = 45147
$ws.Cells.Item(572, 5).Value = 5
$ws.Cells.Item(572, 6).Value = 100112009
$ws.Cells.Item(572, 7).Value = "Acelga"
$ws.Cells.Item(572, 8).Value = "Sin especificar"
$ws.Cells.Item(572, 9).Value = "Primera"
$ws.Cells.Item(572, 10).Value = 120
$ws.Cells.Item(572, 11).Value = 3000
$ws.Cells.Item(572, 12).Value = 3000
$ws.Cells.Item(572, 13).Value = 3000
$ws.Cells.Item(572, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(572, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(572, 16).Value = 500
$ws.Cells.Item(572, 17).Value = 6
$ws.Cells.Item(572, 18).Value = "Hortaliza"
